# Update the "Generate Report for Handback" timestamps.
$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date for the first file (row 2, col G)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-28 05:07:52"

# zh-cn sheet: Correspond Handoff Datetime (col H) / Correspond Handback DateTime (col K) for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-28 05:07:48"
$wsZhCn.Range("K2").Value = "2016-08-28 05:08:22"

# de-de sheet: Correspond Handoff Datetime (col H, shares the same underlying
# shared-string slot as Overview!G2 for this file, so it moves in lockstep)
# and Correspond Handback DateTime (col K) for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-28 05:07:52"
$wsDeDe.Range("K2").Value = "2016-08-28 05:08:28"
